$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values per row
# following re-run of SGNN dialog act annotation.
$ws.Cells.Item(9, 9).Value = "sd"
$ws.Cells.Item(9, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(21, 9).Value = "aa"
$ws.Cells.Item(21, 10).Value = "Agree/Accept"
$ws.Cells.Item(39, 9).Value = "b"
$ws.Cells.Item(39, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(50, 9).Value = "aa"
$ws.Cells.Item(50, 10).Value = "Agree/Accept"
$ws.Cells.Item(60, 9).Value = "b"
$ws.Cells.Item(60, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(67, 9).Value = "sd"
$ws.Cells.Item(67, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(77, 9).Value = "sd"
$ws.Cells.Item(77, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(92, 9).Value = "%"
$ws.Cells.Item(92, 10).Value = "Uninterpretable"
$ws.Cells.Item(94, 9).Value = "%"
$ws.Cells.Item(94, 10).Value = "Uninterpretable"
$ws.Cells.Item(96, 9).Value = "%"
$ws.Cells.Item(96, 10).Value = "Uninterpretable"
$ws.Cells.Item(102, 9).Value = "%"
$ws.Cells.Item(102, 10).Value = "Uninterpretable"
$ws.Cells.Item(109, 9).Value = "sd"
$ws.Cells.Item(109, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(134, 9).Value = "sv"
$ws.Cells.Item(134, 10).Value = "Statement-opinion"
$ws.Cells.Item(144, 9).Value = "aa"
$ws.Cells.Item(144, 10).Value = "Agree/Accept"
$ws.Cells.Item(149, 9).Value = "sd"
$ws.Cells.Item(149, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(155, 9).Value = "sd"
$ws.Cells.Item(155, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(158, 9).Value = "sd"
$ws.Cells.Item(158, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(163, 9).Value = "sv"
$ws.Cells.Item(163, 10).Value = "Statement-opinion"
$ws.Cells.Item(175, 9).Value = "sv"
$ws.Cells.Item(175, 10).Value = "Statement-opinion"
$ws.Cells.Item(178, 9).Value = "ba"
$ws.Cells.Item(178, 10).Value = "Appreciation"
$ws.Cells.Item(180, 9).Value = "sd"
$ws.Cells.Item(180, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(183, 9).Value = "ba"
$ws.Cells.Item(183, 10).Value = "Appreciation"
$ws.Cells.Item(188, 9).Value = "sd"
$ws.Cells.Item(188, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(212, 9).Value = "sv"
$ws.Cells.Item(212, 10).Value = "Statement-opinion"
$ws.Cells.Item(234, 9).Value = "sd"
$ws.Cells.Item(234, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(237, 9).Value = "sv"
$ws.Cells.Item(237, 10).Value = "Statement-opinion"
$ws.Cells.Item(248, 9).Value = "%"
$ws.Cells.Item(248, 10).Value = "Uninterpretable"
$ws.Cells.Item(258, 9).Value = "sd"
$ws.Cells.Item(258, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(260, 9).Value = "sd"
$ws.Cells.Item(260, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(261, 9).Value = "sd"
$ws.Cells.Item(261, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(266, 9).Value = "b"
$ws.Cells.Item(266, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(274, 9).Value = "b"
$ws.Cells.Item(274, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(284, 9).Value = "sv"
$ws.Cells.Item(284, 10).Value = "Statement-opinion"
$ws.Cells.Item(292, 9).Value = "%"
$ws.Cells.Item(292, 10).Value = "Uninterpretable"
$ws.Cells.Item(299, 9).Value = "aa"
$ws.Cells.Item(299, 10).Value = "Agree/Accept"
$ws.Cells.Item(305, 9).Value = "%"
$ws.Cells.Item(305, 10).Value = "Uninterpretable"
$ws.Cells.Item(306, 9).Value = "%"
$ws.Cells.Item(306, 10).Value = "Uninterpretable"
$ws.Cells.Item(320, 9).Value = "sd"
$ws.Cells.Item(320, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(345, 9).Value = "sv"
$ws.Cells.Item(345, 10).Value = "Statement-opinion"
$ws.Cells.Item(354, 9).Value = "aa"
$ws.Cells.Item(354, 10).Value = "Agree/Accept"
$ws.Cells.Item(374, 9).Value = "sd"
$ws.Cells.Item(374, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(384, 9).Value = "sd"
$ws.Cells.Item(384, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(389, 9).Value = "b"
$ws.Cells.Item(389, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(397, 9).Value = "%"
$ws.Cells.Item(397, 10).Value = "Uninterpretable"
$ws.Cells.Item(422, 9).Value = "sd"
$ws.Cells.Item(422, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(423, 9).Value = "b"
$ws.Cells.Item(423, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(431, 9).Value = "sv"
$ws.Cells.Item(431, 10).Value = "Statement-opinion"
$ws.Cells.Item(449, 9).Value = "aa"
$ws.Cells.Item(449, 10).Value = "Agree/Accept"
$ws.Cells.Item(490, 9).Value = "sd"
$ws.Cells.Item(490, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(495, 9).Value = "sv"
$ws.Cells.Item(495, 10).Value = "Statement-opinion"
$ws.Cells.Item(496, 9).Value = "sd"
$ws.Cells.Item(496, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(501, 9).Value = "sd"
$ws.Cells.Item(501, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(502, 9).Value = "ba"
$ws.Cells.Item(502, 10).Value = "Appreciation"
$ws.Cells.Item(509, 9).Value = "sv"
$ws.Cells.Item(509, 10).Value = "Statement-opinion"
$ws.Cells.Item(512, 9).Value = "sv"
$ws.Cells.Item(512, 10).Value = "Statement-opinion"
$ws.Cells.Item(519, 9).Value = "b"
$ws.Cells.Item(519, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(522, 9).Value = "sd"
$ws.Cells.Item(522, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(524, 9).Value = "sd"
$ws.Cells.Item(524, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(529, 9).Value = "sv"
$ws.Cells.Item(529, 10).Value = "Statement-opinion"
$ws.Cells.Item(532, 9).Value = "sv"
$ws.Cells.Item(532, 10).Value = "Statement-opinion"
$ws.Cells.Item(533, 9).Value = "%"
$ws.Cells.Item(533, 10).Value = "Uninterpretable"
$ws.Cells.Item(534, 9).Value = "%"
$ws.Cells.Item(534, 10).Value = "Uninterpretable"
$ws.Cells.Item(535, 9).Value = "sd"
$ws.Cells.Item(535, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(536, 9).Value = "sd"
$ws.Cells.Item(536, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(538, 9).Value = "sd"
$ws.Cells.Item(538, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(541, 9).Value = "sd"
$ws.Cells.Item(541, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(543, 9).Value = "sd"
$ws.Cells.Item(543, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(549, 9).Value = "ba"
$ws.Cells.Item(549, 10).Value = "Appreciation"
$ws.Cells.Item(550, 9).Value = "aa"
$ws.Cells.Item(550, 10).Value = "Agree/Accept"
$ws.Cells.Item(554, 9).Value = "sd"
$ws.Cells.Item(554, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(559, 9).Value = "sd"
$ws.Cells.Item(559, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(563, 9).Value = "aa"
$ws.Cells.Item(563, 10).Value = "Agree/Accept"
$ws.Cells.Item(570, 9).Value = "sd"
$ws.Cells.Item(570, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(572, 9).Value = "sd"
$ws.Cells.Item(572, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(573, 9).Value = "sd"
$ws.Cells.Item(573, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(582, 9).Value = "b"
$ws.Cells.Item(582, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(583, 9).Value = "sd"
$ws.Cells.Item(583, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(600, 9).Value = "%"
$ws.Cells.Item(600, 10).Value = "Uninterpretable"
$ws.Cells.Item(602, 9).Value = "sd"
$ws.Cells.Item(602, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(618, 9).Value = "aa"
$ws.Cells.Item(618, 10).Value = "Agree/Accept"
$ws.Cells.Item(619, 9).Value = "aa"
$ws.Cells.Item(619, 10).Value = "Agree/Accept"
$ws.Cells.Item(621, 9).Value = "sv"
$ws.Cells.Item(621, 10).Value = "Statement-opinion"
$ws.Cells.Item(624, 9).Value = "aa"
$ws.Cells.Item(624, 10).Value = "Agree/Accept"
$ws.Cells.Item(629, 9).Value = "sv"
$ws.Cells.Item(629, 10).Value = "Statement-opinion"
$ws.Cells.Item(644, 9).Value = "sd"
$ws.Cells.Item(644, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(689, 9).Value = "aa"
$ws.Cells.Item(689, 10).Value = "Agree/Accept"
$ws.Cells.Item(692, 9).Value = "aa"
$ws.Cells.Item(692, 10).Value = "Agree/Accept"
$ws.Cells.Item(713, 9).Value = "sd"
$ws.Cells.Item(713, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(723, 9).Value = "sd"
$ws.Cells.Item(723, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(725, 9).Value = "sv"
$ws.Cells.Item(725, 10).Value = "Statement-opinion"
$ws.Cells.Item(729, 9).Value = "b"
$ws.Cells.Item(729, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(747, 9).Value = "sv"
$ws.Cells.Item(747, 10).Value = "Statement-opinion"
$ws.Cells.Item(753, 9).Value = "aa"
$ws.Cells.Item(753, 10).Value = "Agree/Accept"
$ws.Cells.Item(758, 9).Value = "sv"
$ws.Cells.Item(758, 10).Value = "Statement-opinion"
$ws.Cells.Item(764, 9).Value = "sd"
$ws.Cells.Item(764, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(778, 9).Value = "sd"
$ws.Cells.Item(778, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(781, 9).Value = "ba"
$ws.Cells.Item(781, 10).Value = "Appreciation"
$ws.Cells.Item(784, 9).Value = "aa"
$ws.Cells.Item(784, 10).Value = "Agree/Accept"
$ws.Cells.Item(795, 9).Value = "sv"
$ws.Cells.Item(795, 10).Value = "Statement-opinion"
